# Enhance the "课程表" (class schedule) table:
#  - grow the table from 3 columns / 4 rows to 4 columns / 11 rows
#  - re-label the weekday headers and add a 4th weekday column
#  - re-distribute the lesson text/shading across the new layout
#  - center every cell and use the SimSun font for every run

function Get-RgbColor($hex) {
  $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
  $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
  $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
  return $r + ($g * 256) + ($b * 65536)
}

$d = $word.ActiveDocument
$t = $d.Tables(1)
$vt = [char]11   # soft line break (w:br) inside a Range.Text assignment

# ---------------------------------------------------------------------------
# 1) Grow the grid: 3 -> 4 columns, 4 -> 11 rows (new ones appended at the end)
# ---------------------------------------------------------------------------
$t.Columns.Add() | Out-Null

for ($i = 0; $i -lt 7; $i++) {
  $t.Rows.Add() | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Column widths: every column becomes 2160 dxa (108 pt)
# ---------------------------------------------------------------------------
for ($c = 1; $c -le $t.Columns.Count; $c++) {
  $t.Columns($c).Width = 108
}

# ---------------------------------------------------------------------------
# 3) Header row text
# ---------------------------------------------------------------------------
$t.Cell(1,1).Range.Text = "节次"
$t.Cell(1,2).Range.Text = "星期一"
$t.Cell(1,3).Range.Text = "星期二"
$t.Cell(1,4).Range.Text = "星期日"

# ---------------------------------------------------------------------------
# 4) Body rows: period label + the 3 weekday cells (text, shading)
# ---------------------------------------------------------------------------
$periods = @("1","2","3","4","5","6","7","8","9","10")

$mon = @("语文", "", "", "", "", "", "", "", "", "")
$monShd = @("ccffcc","ffffff","ffffff","ffffff","ffffff","ffffff","ffffff","ffffff","ffffff","ffffff")

$tue = @("", ("数学" + $vt + "教师：王老师"), "美术", "", "", "", "", "", "", "数学晚自习")
$tueShd = @("ccffff","ccffff","ffb388","ffffff","ffffff","ffffff","ffffff","ffffff","ffffff","ffe4c4")

$sun = @("数学", "", "", "", "", "", "", "", "", "")
$sunShd = @("ccffff","ffffff","ffffff","ffffff","ffffff","ffffff","ffffff","ffffff","ffffff","ffffff")

for ($i = 0; $i -lt $periods.Count; $i++) {
  $row = $i + 2

  $t.Cell($row,1).Range.Text = $periods[$i]

  $t.Cell($row,2).Range.Text = $mon[$i]
  $t.Cell($row,2).Shading.BackgroundPatternColor = Get-RgbColor $monShd[$i]

  $t.Cell($row,3).Range.Text = $tue[$i]
  $t.Cell($row,3).Shading.BackgroundPatternColor = Get-RgbColor $tueShd[$i]

  $t.Cell($row,4).Range.Text = $sun[$i]
  $t.Cell($row,4).Shading.BackgroundPatternColor = Get-RgbColor $sunShd[$i]
}

# ---------------------------------------------------------------------------
# 5) Whole-table formatting: center every paragraph, SimSun every run
# ---------------------------------------------------------------------------
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell($r,$c)
    $rng = $cell.Range
    $rng.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

    if ($rng.Text.Length -gt 1) {
      $textRng = $d.Range($rng.Start, $rng.End - 1)
      $textRng.Font.Name = "SimSun"
    }
  }
}

Write-Output "done"
